$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 ("age") and row 28 ("education") headers gain a "by ..." /
# "по ..." framing in all three languages (Kyrgyz / Russian / English).
# Written in this order so the shared-string table is appended in the
# same sequence as the target workbook.
$ws.Range("C18").Value = "By age (in years) "
$ws.Range("C28").Value = "By education"
$ws.Range("B18").Value = "По возрасту (в годах)"
$ws.Range("A18").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A28").Value = "Билими боюнча"
$ws.Range("B28").Value = "По образованию"
